$d = $word.ActiveDocument

$replacements = @(
    @("60×67=4020", "58×59=3422"),
    @("86×91=7826", "44×65=2860"),
    @("74×40=2960", "50×53=2650"),
    @("96×81=7776", "35×65=2275"),
    @("92×11=1012", "41×32=1312"),
    @("61×91=5551", "77×34=2618"),
    @("94×40=3760", "17×41=697"),
    @("50×11=550", "99×76=7524"),
    @("73×61=4453", "90×76=6840"),
    @("49×55=2695", "67×94=6298"),
    @("44×15=660", "56×54=3024"),
    @("31×46=1426", "50×86=4300"),
    @("74×78=5772", "11×59=649"),
    @("31×71=2201", "59×50=2950"),
    @("74×35=2590", "30×60=1800"),
    @("49×92=4508", "42×53=2226"),
    @("44×91=4004", "70×23=1610"),
    @("88×29=2552", "51×50=2550"),
    @("87×35=3045", "41×34=1394"),
    @("64×90=5760", "87×89=7743"),
    @("34×71=2414", "50×79=3950"),
    @("44×94=4136", "22×85=1870"),
    @("32×93=2976", "28×30=840"),
    @("98×30=2940", "21×15=315"),
    @("84×93=7812", "95×91=8645")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done: applied $($replacements.Count) replacements"
